$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3631236.8
$ws.Range("J17").Value = 3788842.8
$ws.Range("L17").Value = 11366528.4
$ws.Range("N17").Value = -11366864.4

$ws.Range("H18").Value = 293.125
$ws.Range("I18").Value = 293.125
$ws.Range("K18").Value = 293.125
$ws.Range("M18").Value = -9.125

$ws.Range("H20").Value = 50250
$ws.Range("J20").Value = 39500
$ws.Range("L20").Value = 39500
$ws.Range("N20").Value = -39960

$ws.Range("H28").Value = 602.9231
$ws.Range("J28").Value = 2430.5
$ws.Range("L28").Value = 2430.5
$ws.Range("N28").Value = -3400.5

$ws.Range("H35").Value = 50250
$ws.Range("J35").Value = 39500
$ws.Range("L35").Value = 39500
$ws.Range("N35").Value = -40258

$ws.Range("H38").Value = 732.3182
$ws.Range("I38").Value = 225.55
$ws.Range("J38").Value = 5800
$ws.Range("K38").Value = 676.6500000000001
$ws.Range("L38").Value = 17400
$ws.Range("M38").Value = -304.6500000000001
$ws.Range("N38").Value = -18144

$ws.Range("H63").Value = 49998
$ws.Range("I63").Value = 49998
$ws.Range("K63").Value = 49998
$ws.Range("M63").Value = -49374

$ws.Range("H64").Value = 4248.2
$ws.Range("I64").Value = 3782.6667
$ws.Range("J64").Value = 4447.7144
$ws.Range("K64").Value = 3782.6667
$ws.Range("L64").Value = 4447.7144
$ws.Range("M64").Value = -3534.6667
$ws.Range("N64").Value = -4943.7144

$ws.Range("H66").Value = 49998
$ws.Range("I66").Value = 49998
$ws.Range("K66").Value = 149994
$ws.Range("M66").Value = -146874

$ws.Range("H67").Value = 4248.2
$ws.Range("I67").Value = 3782.6667
$ws.Range("J67").Value = 4447.7144
$ws.Range("K67").Value = 3782.6667
$ws.Range("L67").Value = 4447.7144
$ws.Range("M67").Value = -2924.6667
$ws.Range("N67").Value = -6163.7144

$ws.Range("H76").Value = 3555
$ws.Range("I76").Value = 3604.8333
$ws.Range("J76").Value = 3455.3333
$ws.Range("K76").Value = 3604.8333
$ws.Range("L76").Value = 3455.3333
$ws.Range("M76").Value = -3289.8333
$ws.Range("N76").Value = -4085.3333

$ws.Range("H79").Value = 3555
$ws.Range("I79").Value = 3604.8333
$ws.Range("J79").Value = 3455.3333
$ws.Range("K79").Value = 3604.8333
$ws.Range("L79").Value = 3455.3333
$ws.Range("M79").Value = -2512.8333
$ws.Range("N79").Value = -5639.3333

$ws.Range("H80").Value = 385967.66
$ws.Range("I80").Value = 802.5
$ws.Range("K80").Value = 2407.5
$ws.Range("M80").Value = -1409.5

$ws.Range("H83").Value = 385967.66
$ws.Range("I83").Value = 802.5
$ws.Range("K83").Value = 7222.5
$ws.Range("M83").Value = -2230.5

$ws.Range("H88").Value = 2475.375
$ws.Range("J88").Value = 2274.25
$ws.Range("L88").Value = 2274.25
$ws.Range("N88").Value = -3086.25

$ws.Range("H91").Value = 2475.375
$ws.Range("J91").Value = 2274.25
$ws.Range("L91").Value = 2274.25
$ws.Range("N91").Value = -5082.25

$ws.Range("H107").Value = 4045.1562
$ws.Range("I107").Value = 1084.5834
$ws.Range("J107").Value = 12926.875
$ws.Range("K107").Value = 1084.5834
$ws.Range("L107").Value = 12926.875
$ws.Range("M107").Value = 835.4166
$ws.Range("N107").Value = -16766.875

$ws.Range("H111").Value = 1289.4445
$ws.Range("I111").Value = 1337.4706
$ws.Range("K111").Value = 4012.4118
$ws.Range("M111").Value = -945.4118000000003

$ws.Range("H113").Value = 5482.625
$ws.Range("I113").Value = 4373.4
$ws.Range("J113").Value = 7331.3335
$ws.Range("K113").Value = 4373.4
$ws.Range("L113").Value = 7331.3335
$ws.Range("M113").Value = -1119.4
$ws.Range("N113").Value = -13839.3335

$ws.Range("H141").Value = 3066.9092
$ws.Range("I141").Value = 2733.7144
$ws.Range("K141").Value = 8201.143199999999
$ws.Range("M141").Value = -3021.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3200.8708
$ws.Range("I2").Value = 2594.2632
$ws.Range("K2").Value = 2594.2632
$ws.Range("M2").Value = -2481.2632

$ws.Range("H18").Value = 48124.75

$ws.Range("H32").Value = 5520.213
$ws.Range("I32").Value = 2849.691
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 2849.691
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -2562.691
$ws.Range("N32").Value = -30574

$ws.Range("H74").Value = 10393.075
$ws.Range("I74").Value = 10934.1455
$ws.Range("K74").Value = 10934.1455
$ws.Range("M74").Value = -10060.1455

$ws.Range("H77").Value = 10393.075
$ws.Range("I77").Value = 10934.1455
$ws.Range("K77").Value = 54670.7275
$ws.Range("M77").Value = -50302.7275

$ws.Range("H97").Value = 1036.037
$ws.Range("I97").Value = 1018.52
$ws.Range("J97").Value = 1255
$ws.Range("K97").Value = 1018.52
$ws.Range("L97").Value = 1255
$ws.Range("M97").Value = -522.52
$ws.Range("N97").Value = -2247

$ws.Range("H102").Value = 2182.4614
$ws.Range("I102").Value = 2427.3
$ws.Range("K102").Value = 2427.3
$ws.Range("M102").Value = -805.3000000000002

$ws.Range("H116").Value = 3200.8708
$ws.Range("I116").Value = 2594.2632
$ws.Range("K116").Value = 2594.2632
$ws.Range("M116").Value = -300.2631999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3200.8708
$ws.Range("I3").Value = 2594.2632
$ws.Range("K3").Value = 2594.2632
$ws.Range("M3").Value = -2480.2632

$ws.Range("H20").Value = 15850.087
$ws.Range("I20").Value = 21738.334
$ws.Range("J20").Value = 4809.625
$ws.Range("K20").Value = 21738.334
$ws.Range("L20").Value = 4809.625
$ws.Range("M20").Value = -21491.334
$ws.Range("N20").Value = -5303.625

$ws.Range("H86").Value = 210059.27
$ws.Range("I86").Value = 345520.2
$ws.Range("J86").Value = 3303.158
$ws.Range("K86").Value = 345520.2
$ws.Range("L86").Value = 3303.158
$ws.Range("M86").Value = -344397.2
$ws.Range("N86").Value = -5549.157999999999

$ws.Range("H89").Value = 210059.27
$ws.Range("I89").Value = 345520.2
$ws.Range("J89").Value = 3303.158
$ws.Range("K89").Value = 1727601
$ws.Range("L89").Value = 16515.79
$ws.Range("M89").Value = -1721985
$ws.Range("N89").Value = -27747.79

$ws.Range("H99").Value = 2420.4482
$ws.Range("I99").Value = 2607.28
$ws.Range("J99").Value = 1252.75
$ws.Range("K99").Value = 2607.28
$ws.Range("L99").Value = 1252.75
$ws.Range("M99").Value = -1109.28
$ws.Range("N99").Value = -4248.75

$ws.Range("H105").Value = 4328.28
$ws.Range("I105").Value = 3817.647
$ws.Range("J105").Value = 5413.375
$ws.Range("K105").Value = 3817.647
$ws.Range("L105").Value = 5413.375
$ws.Range("M105").Value = -2070.647
$ws.Range("N105").Value = -8907.375

$ws.Range("H107").Value = 924
$ws.Range("I107").Value = 388.0625
$ws.Range("K107").Value = 388.0625
$ws.Range("M107").Value = 1531.9375

$ws.Range("H132").Value = 77860.75
$ws.Range("J132").Value = 77860.75
$ws.Range("L132").Value = 77860.75
$ws.Range("N132").Value = -87980.75

$ws.Range("H134").Value = 5082.3203
$ws.Range("I134").Value = 2009.295
$ws.Range("K134").Value = 6027.885
$ws.Range("M134").Value = -3492.885

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3151.4
$ws.Range("I16").Value = 3335.5334
$ws.Range("K16").Value = 3335.5334
$ws.Range("M16").Value = -3048.5334

$ws.Range("H31").Value = 2563.8948
$ws.Range("I31").Value = 1228.7778
$ws.Range("K31").Value = 1228.7778
$ws.Range("M31").Value = -933.7778000000001

$ws.Range("H34").Value = 2563.8948
$ws.Range("I34").Value = 1228.7778
$ws.Range("K34").Value = 1228.7778
$ws.Range("M34").Value = -1026.7778

$ws.Range("H50").Value = 47744.75
$ws.Range("J50").Value = 47744.75
$ws.Range("L50").Value = 47744.75
$ws.Range("N50").Value = -48994.75

$ws.Range("H51").Value = 44731.09
$ws.Range("I51").Value = 38545
$ws.Range("J51").Value = 46105.777
$ws.Range("K51").Value = 38545
$ws.Range("L51").Value = 46105.777
$ws.Range("M51").Value = -37809
$ws.Range("N51").Value = -47577.777

$ws.Range("H60").Value = 25277.818
$ws.Range("J60").Value = 35138.57
$ws.Range("L60").Value = 35138.57
$ws.Range("N60").Value = -36160.57

$ws.Range("H61").Value = 44731.09
$ws.Range("I61").Value = 38545
$ws.Range("J61").Value = 46105.777
$ws.Range("K61").Value = 38545
$ws.Range("L61").Value = 46105.777
$ws.Range("M61").Value = -38197
$ws.Range("N61").Value = -46801.777

$ws.Range("H62").Value = 52174.715
$ws.Range("I62").Value = 127493.5
$ws.Range("J62").Value = 5824.6924
$ws.Range("K62").Value = 127493.5
$ws.Range("L62").Value = 5824.6924
$ws.Range("M62").Value = -126869.5
$ws.Range("N62").Value = -7072.6924

$ws.Range("H65").Value = 52174.715
$ws.Range("I65").Value = 127493.5
$ws.Range("J65").Value = 5824.6924
$ws.Range("K65").Value = 637467.5
$ws.Range("L65").Value = 29123.462
$ws.Range("M65").Value = -634347.5
$ws.Range("N65").Value = -35363.462

$ws.Range("H105").Value = 1641.1052
$ws.Range("I105").Value = 1576.9375
$ws.Range("K105").Value = 1576.9375
$ws.Range("M105").Value = 170.0625

$ws.Range("H113").Value = 3151.4
$ws.Range("I113").Value = 3335.5334
$ws.Range("K113").Value = 3335.5334
$ws.Range("M113").Value = -1165.5334

$ws.Range("H132").Value = 20042.145
$ws.Range("I132").Value = 12597.825
$ws.Range("K132").Value = 37793.47500000001
$ws.Range("M132").Value = -35263.47500000001

$ws.Range("H134").Value = 4130.41
$ws.Range("I134").Value = 2532.1177
$ws.Range("K134").Value = 7596.353099999999
$ws.Range("M134").Value = -5061.353099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1973.9688
$ws.Range("J5").Value = 2422.6365
$ws.Range("L5").Value = 7267.9095
$ws.Range("N5").Value = -7491.9095

$ws.Range("H34").Value = 3049.261
$ws.Range("J34").Value = 5707.1665
$ws.Range("L34").Value = 17121.4995
$ws.Range("N34").Value = -17289.4995

$ws.Range("H39").Value = 4627
$ws.Range("J39").Value = 4923.7617
$ws.Range("L39").Value = 14771.2851
$ws.Range("N39").Value = -15359.2851

$ws.Range("H64").Value = 3000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 3000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H117").Value = 1124.9
$ws.Range("I117").Value = 1675
$ws.Range("J117").Value = 987.375
$ws.Range("K117").Value = 5025
$ws.Range("L117").Value = 2962.125
$ws.Range("M117").Value = -1583
$ws.Range("N117").Value = -9846.125

$ws.Range("H135").Value = 1973.9688
$ws.Range("J135").Value = 2422.6365
$ws.Range("L135").Value = 21803.7285
$ws.Range("N135").Value = -26873.7285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 24.333334
$ws.Range("J2").Value = 91.333336
$ws.Range("K2").Value = 24.333334
$ws.Range("L2").Value = 91.333336
$ws.Range("M2").Value = 88.66666599999999
$ws.Range("N2").Value = -317.333336

$ws.Range("H5").Value = 25997.5
$ws.Range("J5").Value = 25998
$ws.Range("L5").Value = 25998
$ws.Range("N5").Value = -26222

$ws.Range("H80").Value = 1849.8572
$ws.Range("I80").Value = 1731.6666
$ws.Range("J80").Value = 1938.5
$ws.Range("K80").Value = 1731.6666
$ws.Range("L80").Value = 1938.5
$ws.Range("M80").Value = -733.6666
$ws.Range("N80").Value = -3934.5

$ws.Range("H83").Value = 1849.8572
$ws.Range("I83").Value = 1731.6666
$ws.Range("J83").Value = 1938.5
$ws.Range("K83").Value = 8658.333000000001
$ws.Range("L83").Value = 9692.5
$ws.Range("M83").Value = -3666.333000000001
$ws.Range("N83").Value = -19676.5

$ws.Range("H122").Value = 1329.5714
$ws.Range("I122").Value = 1237.6364
$ws.Range("K122").Value = 3712.9092
$ws.Range("M122").Value = -1262.9092

$ws.Range("H126").Value = 2844.238
$ws.Range("I126").Value = 2478.1765
$ws.Range("J126").Value = 4400
$ws.Range("K126").Value = 7434.529500000001
$ws.Range("L126").Value = 13200
$ws.Range("M126").Value = -4964.529500000001
$ws.Range("N126").Value = -18140

$ws.Range("H132").Value = 17682.25
$ws.Range("I132").Value = 11378.375
$ws.Range("K132").Value = 34135.125
$ws.Range("M132").Value = -31605.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0

$ws.Range("H7").Value = 4961.846
$ws.Range("I7").Value = 4358.5
$ws.Range("K7").Value = 4358.5
$ws.Range("M7").Value = -4246.5

$ws.Range("H34").Value = 6000
$ws.Range("I34").Value = 6000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5828

$ws.Range("H40").Value = 4199.5
$ws.Range("I40").Value = 3879.4
$ws.Range("K40").Value = 3879.4
$ws.Range("M40").Value = -3743.4

$ws.Range("H68").Value = 2505.2856
$ws.Range("I68").Value = 2997.25
$ws.Range("J68").Value = 1849.3334
$ws.Range("K68").Value = 2997.25
$ws.Range("L68").Value = 1849.3334
$ws.Range("M68").Value = -2248.25
$ws.Range("N68").Value = -3347.3334

$ws.Range("H71").Value = 2505.2856
$ws.Range("I71").Value = 2997.25
$ws.Range("J71").Value = 1849.3334
$ws.Range("K71").Value = 14986.25
$ws.Range("L71").Value = 9246.666999999999
$ws.Range("M71").Value = -11242.25
$ws.Range("N71").Value = -16734.667

$ws.Range("H96").Value = 67999
$ws.Range("J96").Value = 67999
$ws.Range("L96").Value = 67999
$ws.Range("N96").Value = -73491

$ws.Range("H100").Value = 4310.222
$ws.Range("I100").Value = 2984.5715
$ws.Range("K100").Value = 2984.5715
$ws.Range("M100").Value = -2443.5715

$ws.Range("H122").Value = 3509.28
$ws.Range("I122").Value = 3734.077
$ws.Range("K122").Value = 11202.231
$ws.Range("M122").Value = -8752.231

$ws.Range("H126").Value = 4961.846
$ws.Range("I126").Value = 4358.5
$ws.Range("K126").Value = 13075.5
$ws.Range("M126").Value = -10605.5

$ws.Range("H132").Value = 5646.943
$ws.Range("J132").Value = 6858.6665
$ws.Range("L132").Value = 20575.9995
$ws.Range("N132").Value = -25635.9995

$ws.Range("H136").Value = 3872.566
$ws.Range("I136").Value = 3385.5789
$ws.Range("J136").Value = 5106.2666
$ws.Range("K136").Value = 10156.7367
$ws.Range("L136").Value = 15318.7998
$ws.Range("M136").Value = -7606.736699999999
$ws.Range("N136").Value = -20418.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2096.3845
$ws.Range("J81").Value = 2277.5454
$ws.Range("L81").Value = 4555.0908
$ws.Range("N81").Value = -6677.0908

$ws.Range("H84").Value = 2096.3845
$ws.Range("J84").Value = 2277.5454
$ws.Range("L84").Value = 22775.454
$ws.Range("N84").Value = -33383.454

$ws.Range("H107").Value = 1264.45
$ws.Range("I107").Value = 911.7143
$ws.Range("J107").Value = 2087.5
$ws.Range("K107").Value = 2735.1429
$ws.Range("L107").Value = 6262.5
$ws.Range("M107").Value = -815.1428999999998
$ws.Range("N107").Value = -10102.5

$ws.Range("H126").Value = 11051.129
$ws.Range("I126").Value = 5939.391
$ws.Range("K126").Value = 17818.173
$ws.Range("M126").Value = -15348.173

$ws.Range("H132").Value = 156452.61
$ws.Range("I132").Value = 230677.05
$ws.Range("K132").Value = 692031.1499999999
$ws.Range("M132").Value = -689501.1499999999
